$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header labels (B1:D1)
$ws.Range("B1").Value = "train_mae"
$ws.Range("C1").Value = "train_mape"
$ws.Range("D1").Value = "train_rmse"

# Add new header labels (E1:G1), copying the style of the existing headers
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "test_mae"
$ws.Range("F1").Value = "test_mape"
$ws.Range("G1").Value = "test_rmse"

# Row 2 - "Suhu Udara"
$ws.Range("A2").Value = "Suhu Udara (°C)"
$ws.Range("B2").Value = 1.28
$ws.Range("C2").Value = 0.05
$ws.Range("D2").Value = 1.67
$ws.Range("E2").Value = 0.79
$ws.Range("F2").Value = 0.03
$ws.Range("G2").Value = 1.03

# Row 3 - "Kelembapan"
$ws.Range("A3").Value = "Kelembapan (%)"
$ws.Range("B3").Value = 6.84
$ws.Range("C3").Value = 0.07000000000000001
$ws.Range("D3").Value = 9.99
$ws.Range("E3").Value = 3.22
$ws.Range("F3").Value = 0.04
$ws.Range("G3").Value = 5.58
